$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Hunk 1 (row 3)
$ws.Range("H3").Value = 20000
$ws.Range("J3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("N3").Value = -20228
# Hunk 2 (row 48)
$ws.Range("H48").Value = 3000
$ws.Range("J48").Value = 3000
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9584
# Hunk 3 (row 51)
$ws.Range("H51").Value = 2500
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 2500
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 2500
$ws.Range("M51").Value = -2016
$ws.Range("N51").Value = -3468
# Hunk 4 (row 56)
$ws.Range("H56").Value = 3000
$ws.Range("J56").Value = 3000
$ws.Range("L56").Value = 9000
$ws.Range("N56").Value = -10068
# Hunk 5 (row 62)
$ws.Range("H62").Value = 6549.8237
$ws.Range("J62").Value = 7339.2
$ws.Range("L62").Value = 7339.2
$ws.Range("N62").Value = -8587.200000000001
# Hunk 6 (row 65)
$ws.Range("H65").Value = 6549.8237
$ws.Range("J65").Value = 7339.2
$ws.Range("L65").Value = 36696
$ws.Range("N65").Value = -42936
# Hunk 7 (row 86)
$ws.Range("H86").Value = 5631.4546
$ws.Range("I86").Value = 1110
$ws.Range("J86").Value = 15320.286
$ws.Range("K86").Value = 1110
$ws.Range("L86").Value = 15320.286
$ws.Range("M86").Value = 13
$ws.Range("N86").Value = -17566.286
# Hunk 8 (row 89)
$ws.Range("H89").Value = 5631.4546
$ws.Range("I89").Value = 1110
$ws.Range("J89").Value = 15320.286
$ws.Range("K89").Value = 5550
$ws.Range("L89").Value = 76601.42999999999
$ws.Range("M89").Value = 66
$ws.Range("N89").Value = -87833.42999999999
# Hunk 9 (row 100)
$ws.Range("H100").Value = 250003380
$ws.Range("I100").Value = 1000000000
$ws.Range("J100").Value = 4502
$ws.Range("K100").Value = 1000000000
$ws.Range("L100").Value = 4502
$ws.Range("M100").Value = -999999459
$ws.Range("N100").Value = -5584
# Hunk 10 (row 102)
$ws.Range("H102").Value = 20000
$ws.Range("J102").Value = 20000
$ws.Range("L102").Value = 20000
$ws.Range("N102").Value = -26490
# Hunk 11 (row 107)
$ws.Range("H107").Value = 654.4815
$ws.Range("I107").Value = 425
$ws.Range("K107").Value = 425
$ws.Range("M107").Value = 1495
# Hunk 12 (row 116)
$ws.Range("H116").Value = 5958.8335
$ws.Range("J116").Value = 7111.778
$ws.Range("L116").Value = 7111.778
$ws.Range("N116").Value = -13995.778
# Hunk 13 (row 132)
$ws.Range("H132").Value = 1987.4286
$ws.Range("I132").Value = 2117.8667
$ws.Range("J132").Value = 520
$ws.Range("K132").Value = 6353.6001
$ws.Range("L132").Value = 1560
$ws.Range("M132").Value = -3823.6001
$ws.Range("N132").Value = -6620
# Hunk 14 (row 137)
$ws.Range("H137").Value = 2506.1428
$ws.Range("I137").Value = 2196.6667
$ws.Range("J137").Value = 2738.25
$ws.Range("K137").Value = 6590.000100000001
$ws.Range("L137").Value = 8214.75
$ws.Range("M137").Value = -4040.000100000001
$ws.Range("N137").Value = -13314.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Hunk 15 (row 32)
$ws.Range("H32").Value = 8151.6416
$ws.Range("I32").Value = 5557.341
$ws.Range("K32").Value = 5557.341
$ws.Range("M32").Value = -5270.341
# Hunk 16 (row 61)
$ws.Range("H61").Value = 1219.3438
$ws.Range("I61").Value = 1219.5416
$ws.Range("J61").Value = 1218.75
$ws.Range("K61").Value = 1219.5416
$ws.Range("L61").Value = 1218.75
$ws.Range("M61").Value = -1007.5416
$ws.Range("N61").Value = -1642.75
# Hunk 17 (row 103)
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
# Hunk 18 (row 110)
$ws.Range("H110").Value = 666.5
$ws.Range("I110").Value = 627.6667
$ws.Range("J110").Value = 783
$ws.Range("K110").Value = 627.6667
$ws.Range("L110").Value = 783
$ws.Range("M110").Value = 1417.3333
$ws.Range("N110").Value = -4873
# Hunk 19 (row 132)
$ws.Range("H132").Value = 12071.574
$ws.Range("I132").Value = 1156.3
$ws.Range("K132").Value = 3468.9
$ws.Range("M132").Value = -938.8999999999996
# Hunk 20 (row 136)
$ws.Range("H136").Value = 1219.3438
$ws.Range("I136").Value = 1219.5416
$ws.Range("J136").Value = 1218.75
$ws.Range("K136").Value = 3658.6248
$ws.Range("L136").Value = 3656.25
$ws.Range("M136").Value = -1108.6248
$ws.Range("N136").Value = -8756.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Hunk 21 (row 31)
$ws.Range("H31").Value = 3785.4285
$ws.Range("I31").Value = 5387.3335
$ws.Range("K31").Value = 5387.3335
$ws.Range("M31").Value = -5092.3335
# Hunk 22 (row 34)
$ws.Range("H34").Value = 3785.4285
$ws.Range("I34").Value = 5387.3335
$ws.Range("K34").Value = 5387.3335
$ws.Range("M34").Value = -5185.3335
# Hunk 23 (row 50)
$ws.Range("H50").Value = 19975
$ws.Range("J50").Value = 19975
$ws.Range("L50").Value = 19975
$ws.Range("N50").Value = -21225
# Hunk 24 (row 58)
$ws.Range("H58").Value = 18180.334
$ws.Range("J58").Value = 43327.168
$ws.Range("L58").Value = 43327.168
$ws.Range("N58").Value = -43733.168
# Hunk 25 (row 86)
$ws.Range("H86").Value = 22197.889
$ws.Range("I86").Value = 10621.4
$ws.Range("K86").Value = 10621.4
$ws.Range("M86").Value = -9498.4
# Hunk 26 (row 89)
$ws.Range("H89").Value = 22197.889
$ws.Range("I89").Value = 10621.4
$ws.Range("K89").Value = 53107
$ws.Range("M89").Value = -47491
# Hunk 27 (row 132)
$ws.Range("H132").Value = 2458.6365
$ws.Range("I132").Value = 1670.7037
$ws.Range("J132").Value = 6004.3335
$ws.Range("K132").Value = 5012.1111
$ws.Range("L132").Value = 18013.0005
$ws.Range("M132").Value = -2482.1111
$ws.Range("N132").Value = -23073.0005
# Hunk 28 (row 136)
$ws.Range("H136").Value = 18180.334
$ws.Range("J136").Value = 43327.168
$ws.Range("L136").Value = 129981.504
$ws.Range("N136").Value = -135081.504

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Hunk 29 (row 113)
$ws.Range("H113").Value = 597.6667
$ws.Range("I113").Value = 649.25
$ws.Range("J113").Value = 494.5
$ws.Range("K113").Value = 1947.75
$ws.Range("L113").Value = 1483.5
$ws.Range("M113").Value = 222.25
$ws.Range("N113").Value = -5823.5
# Hunk 30 (row 131)
$ws.Range("H131").Value = 730.76
$ws.Range("J131").Value = 732.0808
$ws.Range("L131").Value = 2196.2424
$ws.Range("N131").Value = -12276.2424

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Hunk 31 (row 36)
$ws.Range("H36").Value = 1633.3334
$ws.Range("I36").Value = 1200
$ws.Range("J36").Value = 2500
$ws.Range("K36").Value = 1200
$ws.Range("L36").Value = 2500
$ws.Range("M36").Value = -715
$ws.Range("N36").Value = -3470

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Hunk 32 (row 64)
$ws.Range("H64").Value = 30783.5
$ws.Range("I64").Value = 10136
$ws.Range("J64").Value = 37666
$ws.Range("K64").Value = 10136
$ws.Range("L64").Value = 37666
$ws.Range("M64").Value = -9911
$ws.Range("N64").Value = -38116
# Hunk 33 (row 67)
$ws.Range("H67").Value = 30783.5
$ws.Range("I67").Value = 10136
$ws.Range("J67").Value = 37666
$ws.Range("K67").Value = 10136
$ws.Range("L67").Value = 37666
$ws.Range("M67").Value = -9356
$ws.Range("N67").Value = -39226

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Hunk 34 (row 82)
$ws.Range("H82").Value = 8914
$ws.Range("I82").Value = 2273
$ws.Range("K82").Value = 2273
$ws.Range("M82").Value = -1890
# Hunk 35 (row 85)
$ws.Range("H85").Value = 8914
$ws.Range("I85").Value = 2273
$ws.Range("K85").Value = 2273
$ws.Range("M85").Value = -947
# Hunk 36 (row 136)
$ws.Range("H136").Value = 20835530
$ws.Range("I136").Value = 31251146
$ws.Range("J136").Value = 4294.0625
$ws.Range("K136").Value = 93753438
$ws.Range("L136").Value = 12882.1875
$ws.Range("M136").Value = -93750888
$ws.Range("N136").Value = -17982.1875
